# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect the latest scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 6344
    4  = 179
    7  = 1905
    10 = 963
    11 = 273
    12 = 5591
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
